$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for A2:A86 (row i -> value at index i-2)
$newValues = @(
    35.36355261064645,  # A2
    38.48840163139249,  # A3
    40.33928021310246,  # A4
    41.13250984825467,  # A5
    36.82726021617131,  # A6
    41.43643303789648,  # A7
    38.62649634902951,  # A8
    37.29151186388066,  # A9
    36.38762104345184,  # A10
    36.26878131643321,  # A11
    31.41493417656761,  # A12
    35.03963592425407,  # A13
    33.61008544654723,  # A14
    26.06827386074252,  # A15
    25.6504450875162,  # A16
    23.07929688879923,  # A17
    23.07929688879923,  # A18
    23.07929688879923,  # A19
    23.07929688879923,  # A20
    23.07929688879923,  # A21
    23.07929688879923,  # A22
    23.07929688879923,  # A23
    23.07929688879923,  # A24
    23.07929688879923,  # A25
    23.07929688879923,  # A26
    23.07929688879923,  # A27
    23.07929688879923,  # A28
    23.07929688879923,  # A29
    23.07929688879923,  # A30
    23.07929688879923,  # A31
    23.07929688879923,  # A32
    23.07929688879923,  # A33
    23.07929688879923,  # A34
    23.07929688879923,  # A35
    23.07929688879923,  # A36
    23.07929688879923,  # A37
    23.07929688879923,  # A38
    23.07929688879923,  # A39
    23.07929688879923,  # A40
    23.07929688879923,  # A41
    23.07929688879923,  # A42
    23.07929688879923,  # A43
    23.07929688879923,  # A44
    23.07929688879923,  # A45
    23.07929688879923,  # A46
    23.07929688879923,  # A47
    23.07929688879923,  # A48
    23.07929688879923,  # A49
    23.07929688879923,  # A50
    23.07929688879923,  # A51
    23.07929688879923,  # A52
    23.07929688879923,  # A53
    23.07929688879923,  # A54
    23.07929688879923,  # A55
    23.07929688879923,  # A56
    23.07929688879923,  # A57
    23.07929688879923,  # A58
    23.07929688879923,  # A59
    23.07929688879923,  # A60
    23.07929688879923,  # A61
    23.07929688879923,  # A62
    23.07929688879923,  # A63
    23.07929688879923,  # A64
    23.07929688879923,  # A65
    23.07929688879923,  # A66
    23.07929688879923,  # A67
    23.07929688879923,  # A68
    23.07929688879923,  # A69
    23.07929688879923,  # A70
    23.07929688879923,  # A71
    23.07929688879923,  # A72
    23.07929688879923,  # A73
    23.07929688879923,  # A74
    23.07929688879923,  # A75
    23.07929688879923,  # A76
    23.07929688879923,  # A77
    23.07929688879923,  # A78
    23.07929688879923,  # A79
    23.07929688879923,  # A80
    23.07929688879923,  # A81
    23.07929688879923,  # A82
    23.07929688879923,  # A83
    23.07929688879923,  # A84
    23.07929688879923,  # A85
    23.07929688879923  # A86
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
